# TestOversigt.xlsx - Ark1 opdateret med testnavne, udført-af samt ny
# kolonnebredde/formatering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row content: Test navn (B) and Udført af (D) --------------------------
$testNavne = @{
    2  = "Søg på en eksisterende ordre"
    3  = "Søg på en ordre som ikke eksisterer"
    4  = "Søg på en eksisterende faktura"
    5  = "Søg på en faktura som ikke eksisterer"
    6  = "Test at systemet giver en fejlbesked når der indtastes bogstaver i søgefeltet"
    7  = "Opgrader en eksisterende ordre til en faktura"
    8  = "Opgrader en eksisterende ordre til en faktura og tilføj en bedemand."
    9  = 'Søg på en fakturanummer uden at tilføje "-" i fakturanummeret'
    10 = "Slet en eksisterende ordre"
    11 = "Slet en eksisterende faktura"
    12 = "Test at man kan oprette en ny bedemand "
    13 = "Ændre en faktura til at være betalt"
}

# Rows whose "Test navn" cell uses the special Cambria-font style.
$cambriaRows = @(6, 9, 12, 13)

# 1) Give the Cambria-styled cells their font first, so the new font and the
#    matching cell style land in the workbook before anything else is added.
foreach ($r in $cambriaRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $testNavne[$r]
    $cell.Font.Name = "Cambria"
}

# 2) Fill in the remaining "Test navn" cells with the plain (default) style.
for ($r = 2; $r -le 13; $r++) {
    if ($cambriaRows -contains $r) { continue }
    $ws.Cells.Item($r, 2).Value = $testNavne[$r]
}

# 3) "Udført af" column (D) for rows 2-13 -> "Anette".
for ($r = 2; $r -le 13; $r++) {
    $d = $ws.Cells.Item($r, 4)
    $d.Value = "Anette"
    $d.HorizontalAlignment = -4108
}

# 4) Re-colour column A (Test id) for rows 2-41 with the same fill used by the
#    header row, keeping the normal (non-bold) font and centred alignment.
$ws.Range("A1").Copy()
for ($r = 2; $r -le 41; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.PasteSpecial(-4122)
    $a.Font.Bold = $false
    $a.Font.Size = 11
}
$excel.CutCopyMode = $false

# 5) Column widths: A narrower, B much wider to fit the new test descriptions.
$ws.Columns.Item(1).ColumnWidth = 16.5703125
$ws.Columns.Item(2).ColumnWidth = 75.140625
